# Rename the header row:
#   *_old  -> *_FV2310
#   *_new  -> *_FV2404
# (column K "diff" is left untouched)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @{
  "A1" = "Segmentname_FV2310"
  "B1" = "Segmentgruppe_FV2310"
  "C1" = "Segment_FV2310"
  "D1" = "Datenelement_FV2310"
  "E1" = "Segment ID_FV2310"
  "F1" = "Code_FV2310"
  "G1" = "Qualifier_FV2310"
  "H1" = "Beschreibung_FV2310"
  "I1" = "Bedingungsausdruck_FV2310"
  "J1" = "Bedingung_FV2310"
  "L1" = "Segmentname_FV2404"
  "M1" = "Segmentgruppe_FV2404"
  "N1" = "Segment_FV2404"
  "O1" = "Datenelement_FV2404"
  "P1" = "Segment ID_FV2404"
  "Q1" = "Code_FV2404"
  "R1" = "Qualifier_FV2404"
  "S1" = "Beschreibung_FV2404"
  "T1" = "Bedingungsausdruck_FV2404"
  "U1" = "Bedingung_FV2404"
}
foreach ($addr in $headers.Keys) {
    $ws.Range($addr).Value = $headers[$addr]
}

# Turn the data range into an Excel table ("Table1") so the renamed
# headers are also reflected in the table column definitions.
$dataRange = $ws.Range("A1:U64")
$table = $ws.ListObjects.Add(1, $dataRange, 0, 1)
$table.Name = "Table1"

# Freeze the header row (freeze panes at A2 -> row 1 stays visible).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
